$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 603.6579
$ws.Range("I17").Value = 476.14285
$ws.Range("J17").Value = 632.4516
$ws.Range("K17").Value = 1428.42855
$ws.Range("L17").Value = 1897.3548
$ws.Range("M17").Value = -1260.42855
$ws.Range("N17").Value = -2233.3548

$ws.Range("H40").Value = 1790.5264
$ws.Range("I40").Value = 1623.3334
$ws.Range("J40").Value = 1821.875
$ws.Range("K40").Value = 1623.3334
$ws.Range("L40").Value = 1821.875
$ws.Range("M40").Value = -1448.3334
$ws.Range("N40").Value = -2171.875

$ws.Range("H64").Value = 3345
$ws.Range("I64").Value = 3350
$ws.Range("J64").Value = 3330
$ws.Range("K64").Value = 3350
$ws.Range("L64").Value = 3330
$ws.Range("M64").Value = -3102
$ws.Range("N64").Value = -3826

$ws.Range("H67").Value = 3345
$ws.Range("I67").Value = 3350
$ws.Range("J67").Value = 3330
$ws.Range("K67").Value = 3350
$ws.Range("L67").Value = 3330
$ws.Range("M67").Value = -2492
$ws.Range("N67").Value = -5046

$ws.Range("H98").Value = 1406
$ws.Range("I98").Value = 1496.6666
$ws.Range("J98").Value = 998
$ws.Range("K98").Value = 1496.6666
$ws.Range("L98").Value = 998
$ws.Range("M98").Value = 1.333399999999983
$ws.Range("N98").Value = -3994

$ws.Range("H111").Value = 2937.9
$ws.Range("I111").Value = 5307.25
$ws.Range("J111").Value = 1358.3334
$ws.Range("K111").Value = 15921.75
$ws.Range("L111").Value = 4075.0002
$ws.Range("M111").Value = -12854.75
$ws.Range("N111").Value = -10209.0002

$ws.Range("H122").Value = 1406
$ws.Range("I122").Value = 1496.6666
$ws.Range("J122").Value = 998
$ws.Range("K122").Value = 4489.9998
$ws.Range("L122").Value = 2994
$ws.Range("M122").Value = -2039.9998
$ws.Range("N122").Value = -7894

$ws.Range("H125").Value = 1046.75
$ws.Range("I125").Value = 483
$ws.Range("J125").Value = 1234.6666
$ws.Range("K125").Value = 4347
$ws.Range("L125").Value = 11111.9994
$ws.Range("M125").Value = -1887
$ws.Range("N125").Value = -16031.9994

$ws.Range("H132").Value = 4763650
$ws.Range("I132").Value = 5716199.5
$ws.Range("J132").Value = 903.2
$ws.Range("K132").Value = 17148598.5
$ws.Range("L132").Value = 2709.6
$ws.Range("M132").Value = -17146068.5
$ws.Range("N132").Value = -7769.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1209.5714
$ws.Range("I2").Value = 918
$ws.Range("J2").Value = 5000
$ws.Range("K2").Value = 918
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = -805
$ws.Range("N2").Value = -5226

$ws.Range("H45").Value = 1258.2
$ws.Range("I45").Value = 1379.4
$ws.Range("J45").Value = 1137
$ws.Range("K45").Value = 1379.4
$ws.Range("L45").Value = 1137
$ws.Range("M45").Value = -1002.4
$ws.Range("N45").Value = -1891

$ws.Range("H61").Value = 2166.9473
$ws.Range("I61").Value = 1557.76
$ws.Range("J61").Value = 3338.4614
$ws.Range("K61").Value = 1557.76
$ws.Range("L61").Value = 3338.4614
$ws.Range("M61").Value = -1345.76
$ws.Range("N61").Value = -3762.4614

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H116").Value = 1209.5714
$ws.Range("I116").Value = 918
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 918
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = 1376
$ws.Range("N116").Value = -9588

$ws.Range("H136").Value = 2166.9473
$ws.Range("I136").Value = 1557.76
$ws.Range("J136").Value = 3338.4614
$ws.Range("K136").Value = 4673.28
$ws.Range("L136").Value = 10015.3842
$ws.Range("M136").Value = -2123.28
$ws.Range("N136").Value = -15115.3842

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1209.5714
$ws.Range("I3").Value = 918
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 918
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = -804
$ws.Range("N3").Value = -5228

$ws.Range("H20").Value = 2524.9312
$ws.Range("I20").Value = 2874.3845
$ws.Range("J20").Value = 2241
$ws.Range("K20").Value = 2874.3845
$ws.Range("L20").Value = 2241
$ws.Range("M20").Value = -2627.3845
$ws.Range("N20").Value = -2735

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 861.9167
$ws.Range("I16").Value = 821.6667
$ws.Range("J16").Value = 902.1667
$ws.Range("K16").Value = 821.6667
$ws.Range("L16").Value = 902.1667
$ws.Range("M16").Value = -534.6667
$ws.Range("N16").Value = -1476.1667

$ws.Range("H58").Value = 711.1613
$ws.Range("I58").Value = 656.2727
$ws.Range("J58").Value = 845.3333
$ws.Range("K58").Value = 656.2727
$ws.Range("L58").Value = 845.3333
$ws.Range("M58").Value = -453.2727
$ws.Range("N58").Value = -1251.3333

$ws.Range("H86").Value = 142861500
$ws.Range("I86").Value = 333335170
$ws.Range("K86").Value = 333335170
$ws.Range("M86").Value = -333334047

$ws.Range("H89").Value = 142861500
$ws.Range("I89").Value = 333335170
$ws.Range("K89").Value = 1666675850
$ws.Range("M89").Value = -1666670234

$ws.Range("H99").Value = 2614.7058
$ws.Range("I99").Value = 2440
$ws.Range("J99").Value = 2687.5
$ws.Range("K99").Value = 2440
$ws.Range("L99").Value = 2687.5
$ws.Range("M99").Value = -942
$ws.Range("N99").Value = -5683.5

$ws.Range("H107").Value = 676.1667
$ws.Range("I107").Value = 645.53845
$ws.Range("J107").Value = 712.36365
$ws.Range("K107").Value = 645.53845
$ws.Range("L107").Value = 712.36365
$ws.Range("M107").Value = 1274.46155
$ws.Range("N107").Value = -4552.36365

$ws.Range("H113").Value = 861.9167
$ws.Range("I113").Value = 821.6667
$ws.Range("J113").Value = 902.1667
$ws.Range("K113").Value = 821.6667
$ws.Range("L113").Value = 902.1667
$ws.Range("M113").Value = 1348.3333
$ws.Range("N113").Value = -5242.1667

$ws.Range("H126").Value = 2614.7058
$ws.Range("I126").Value = 2440
$ws.Range("J126").Value = 2687.5
$ws.Range("K126").Value = 7320
$ws.Range("L126").Value = 8062.5
$ws.Range("M126").Value = -4850
$ws.Range("N126").Value = -13002.5

$ws.Range("H134").Value = 1692.3572
$ws.Range("I134").Value = 1479.909
$ws.Range("J134").Value = 2471.3333
$ws.Range("K134").Value = 4439.727000000001
$ws.Range("L134").Value = 7413.999899999999
$ws.Range("M134").Value = -1904.727000000001
$ws.Range("N134").Value = -12483.9999

$ws.Range("H136").Value = 711.1613
$ws.Range("I136").Value = 656.2727
$ws.Range("J136").Value = 845.3333
$ws.Range("K136").Value = 1968.8181
$ws.Range("L136").Value = 2535.9999
$ws.Range("M136").Value = 581.1819
$ws.Range("N136").Value = -7635.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1510.5555
$ws.Range("I92").Value = 600
$ws.Range("J92").Value = 2090
$ws.Range("K92").Value = 1800
$ws.Range("L92").Value = 6270
$ws.Range("M92").Value = -552
$ws.Range("N92").Value = -8766

$ws.Range("H107").Value = 505
$ws.Range("I107").Value = 447.5
$ws.Range("J107").Value = 562.5
$ws.Range("K107").Value = 1342.5
$ws.Range("L107").Value = 1687.5
$ws.Range("M107").Value = 577.5
$ws.Range("N107").Value = -5527.5

$ws.Range("H113").Value = 547.7442
$ws.Range("I113").Value = 532.6667
$ws.Range("J113").Value = 562.13635
$ws.Range("K113").Value = 1598.0001
$ws.Range("L113").Value = 1686.40905
$ws.Range("M113").Value = 571.9999
$ws.Range("N113").Value = -6026.40905

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4062.5
$ws.Range("I126").Value = 4600
$ws.Range("J126").Value = 3525
$ws.Range("K126").Value = 13800
$ws.Range("L126").Value = 10575
$ws.Range("M126").Value = -11330
$ws.Range("N126").Value = -15515

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1400
$ws.Range("I61").Value = 1350
$ws.Range("J61").Value = 1500
$ws.Range("K61").Value = 1350
$ws.Range("L61").Value = 1500
$ws.Range("M61").Value = -1148
$ws.Range("N61").Value = -1904

$ws.Range("H113").Value = 1400
$ws.Range("I113").Value = 1350
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1350
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 820
$ws.Range("N113").Value = -5840

$ws.Range("H132").Value = 1532.3334
$ws.Range("I132").Value = 1315.2667
$ws.Range("J132").Value = 2255.889
$ws.Range("K132").Value = 3945.800099999999
$ws.Range("L132").Value = 6767.667
$ws.Range("M132").Value = -1415.800099999999
$ws.Range("N132").Value = -11827.667
